$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C2').Value = 'RT @BALKIRAGA: PukKodu ATATÜRK
Aydınlık izleri silinmedikçe ülke doğruyu bulacaktır.
@MKAtimi 
@hzlandrc 
@B6tur 
@YildirimUgurgul 
@zzoguz…'
$ws.Range('C3').Value = 'RT @tancabrona: — ¿Sigues enojada?
— No.
— ¿Y ese cuchillo?
— https://t.co/oClipdteUA'
$ws.Range('C4').Value = 'RT @BLACKPINKGLOBAL: [170327] [PRESS] Jennie at BOON THE SHOP Event #BLACKPINK #블랙핑크 #JENNIE #제니 https://t.co/oT8cI2dDFb'
$ws.Range('C5').Value = 'RT @12Super1Hero: つじ写真館さんに昨日お忙しい時間に行きました！
びゅうおの写真とコメントをお願いして帰宅しました
お茶ありがとうございましたm(*_ _)m https://t.co/sDvwhF9zVD'
$ws.Range('C6').Value = 'RT @The40Chambers: CRICE https://t.co/UySmRmjKo7'
$ws.Range('C7').Value = 'What a fucking iconic queen. https://t.co/fApuhGZs0b'
$ws.Range('C8').Value = 'RT @_omanprojects: ??شركة تنمية نفط #عمان تطرح مشروع تطوير رأس الحمراء للاستثمار.. https://t.co/Jid0Plt2zZ'
$ws.Range('C9').Value = 'RT @gblardone: -2 fois ministre de Fillon (+1,2M chômeurs +600MM dette)
-2 fois déjà sous Chirac
-2 fois porte-parole du gouv. dès 1995
-Dé…'
$ws.Range('C10').Value = 'RT @debbiemc1547: https://t.co/ZpnA78zNT3'
$ws.Range('C11').Value = 'RT @AuntyTalks: இந்த மாதிரி ரசிச்சு புண்டைய நக்க யாராவது lesbo Twitter ல இருக்கீங்களா ப்ரண்ட்ஸ்?? https://t.co/lHDcVKw1cB'
$ws.Range('C12').Value = 'RT @kacsaatolldu: Anasınıfına bile Erdoğan''lı kitap dağıtanlar, Nutuk''u ''siyasi propaganda olur''diye yasaklayalı birkaç gün oldu https://t.…'
$ws.Range('C13').Value = 'RT @VoetbalInside: PRIMEUR: Hierbij presenteren wij jullie de nieuwe staff van @OnsOranje / @KNVB. Beter? ?? ???? #voetbalinside https://t.co/…'
$ws.Range('C14').Value = 'https://t.co/FwAYHcSjjx'
$ws.Range('C15').Value = 'RT @n_nammimi: ผัวเดย์ค่ะ แจกเป๋านี้ 1 ใบ รีไปนะ อิอิ #KINGJACKSONDAY https://t.co/lNsQthhW4m'
$ws.Range('C16').Value = 'RT @NTelevisa_com: El premio llega en un momento difícil para México porque han sido asesinados 3 periodista en el último mes:@CarlosLoret…'
$ws.Range('C17').Value = 'RT @Iovekth: this angel https://t.co/gV1dqldcsK'
$ws.Range('C18').Value = 'Gemüse-Implantate: So sollen aus Spinat Ersatzteile für unser Herz wachsen https://t.co/lTHuSSSsfc https://t.co/DQJhWzysVk'
$ws.Range('C19').Value = 'RT @tuanarchives: happy birthday to got7''s angel, he deserves to be happy and loved, we love him so much, he is always good to mark ✨❤ #KIN…'
$ws.Range('C20').Value = 'https://t.co/g1OrhO1Fa1'
$ws.Range('C21').Value = 'They''re in Estonia working for Putin. There''s NO WAY on EARTH Trump''s got 36%, Putinbots are hacking the pollsters… https://t.co/3oZUntzOZA'
$ws.Range('C22').Value = 'RT @jenarovillamil: "Lo contrario de la libertad no es el determinismo sino el fatalismo": Jean Paul Sartre https://t.co/xRHAvr1FcD'
$ws.Range('C23').Value = 'RT @deray: what type of society raised this domestic terrorist? https://t.co/XI501hEBVl'
$ws.Range('C24').Value = 'RT @calzonaflames: "If you love someone, you tell them" #12YearsOfGreysAnatomy *sobs* https://t.co/p7YBmhMBlx'
$ws.Range('C25').Value = 'RT @Dalton_Chad: Come to the library mall and #VoteColeUmeh !! Dex did, you should too. https://t.co/aiWPHoOzV6'
$ws.Range('C26').Value = 'RT @RocketJoy: Check out the insides of our Crew Dragon spacecraft and the system that will support human life in space! https://t.co/xUnmE…'
$ws.Range('C27').Value = 'En el hilo de la semana, pude descargarme está joya. Aprovechen  https://t.co/mmXAYq6QsE https://t.co/lJ0s6I4td1'
$ws.Range('C28').Value = 'RT @hehehe9988: โอ้ยขรรม?? เจ๊หนิงไลฟ์สดของพี่เป๊กอยู่ แล้วซูมไปที่เป้าเพราะพี่เป๊กบอกพึ่งรู้ตัวว่าเป้าแตก55555555555555 #เป๊กผลิตโชค #ผลิต…'
$ws.Range('C29').Value = 'これほしい(´･ω･`)
カラコンも欲しいし。
今更ながらRuuaのカバンよりこっちにすればよかった、、かも https://t.co/Z6ZrEl8rJD'
$ws.Range('C30').Value = 'Please present your evidence that causality holds outside of this universe. https://t.co/v0k8R0sPSK'
$ws.Range('C31').Value = '【モンスト】『わくりん2倍』で金種大量出現！勝ち組プロスト多数爆誕ｷﾀ━━━━(ﾟ∀ﾟ)━━━━!!【画像あり】⇒ https://t.co/WkI6UNUMrl https://t.co/H1rLHWHezq'
$ws.Range('C32').Value = 'https://t.co/NRLkmYTh78'
$ws.Range('C33').Value = 'Skrenggeh! https://t.co/E34n7LOTP1'
$ws.Range('C34').Value = 'RT @Jukeslol: SALVE RAPEIZE
ACABAIE DE ACORDAR E STREAM JA ESTA ONLINE
JOGANDO NA MAIN CHALLENGER
#enois
https://t.co/gm1jefJqs8 https://t.…'
$ws.Range('C35').Value = 'https://t.co/fyFYjKxDqP'
$ws.Range('C36').Value = 'RT @PaosameSurTwi: Bientôt les gens ils porteront ca oklm dans la rue https://t.co/iQhCA14fnL'
$ws.Range('C37').Value = 'RT @SeamusGorman1: @Meme_Druid @ReeReeC2 @RandallKraft @Mom2theCorps @ggentlemanirish @Kimmie091577 @JewelsJones1 @jimmygarner @KeecoWang5…'
$ws.Range('C38').Value = 'RT @atletico: Gostou do clipe #Galo109? Então assista ao making of com os erros de gravação! #Galo #PaixãoDoPovo https://t.co/LzydwyUDF3'
$ws.Range('C39').Value = 'This man is only accused of 5 billion ruppies of corruption.
Still a shoe of Zardari 
 https://t.co/gRDluzDvLq'
$ws.Range('C40').Value = 'RT @WBCBaseball: Now’s your opportunity to win a sweet prize! RT for a chance to win this! #WBC2017 https://t.co/H9S34RYgrU'
$ws.Range('C41').Value = 'RT @jvaldez666: Eres flor eres hermosa ?? https://t.co/HB8npfvnFb'
$ws.Range('C42').Value = 'RT @kirbsterr__: Fr though ?????? I''m fucking tripping https://t.co/8CQXKw4h7G'
$ws.Range('C43').Value = 'RT @ItsFoodPorn: Blueberry Cheesecake https://t.co/fyhNhlmBts'
$ws.Range('C44').Value = 'RT @guillaumecastan: Pour lutter contre la concentration des terres par les grandes sociétés agricoles et aider les jeunes agriculteurs, #J…'
$ws.Range('C45').Value = 'RT @tarrraan_: Come on y''all !! Get me my dream dog ! Please !! https://t.co/nZctc4jIhC'
$ws.Range('C46').Value = 'nahhh lol https://t.co/AOEzaEJaCU'
$ws.Range('C47').Value = '#ابسط_حقوق_المصريين https://t.co/oYx7LpgJUr'
$ws.Range('C48').Value = 'RT @TolgaYakali: ??31.03.2017 https://t.co/Ln2FQFfU9h'
$ws.Range('C49').Value = 'RT @baptista1904: QUAL É O MELHOR CLUBE DO MUNDO https://t.co/6HQ4kYMXGx'
$ws.Range('C50').Value = 'RT @jiminspired: the days where jimin was chubby and buff, the cutest basketball player https://t.co/05yjrLZSXn'
$ws.Range('C51').Value = 'RT @fyo101: اذا اخوي نام في الصالة
 وامي قالت وده سريره ?? https://t.co/KCkYzXMfFJ'
$ws.Range('C52').Value = 'RT @camucha4: URGENTE! Necesitamos Avastin 400 (4 frascos) , para que Melody pueda continuar su tratamiento.  #TodosPorMelody Gracias RT #u…'
$ws.Range('C53').Value = 'RT @KFCBarstool: God I love Frank Martin https://t.co/euQXXutVg3'
$ws.Range('C54').Value = 'RT @nia4_trump: On #MuslimWomensDay let''s take a moment to understand the etiquette of Wife Beatings according to Islam &amp; Sharia.
https://t…'
$ws.Range('C55').Value = 'ENJOY EVERY STEP ALONG THE WAY! https://t.co/ULrzDrwv2D'
$ws.Range('C56').Value = 'RT @SytnerBMW: Check out this M4 Convertible in Sakhir Orange with Black Individual Merino Leather at Sytner Sheffield. For more info, call…'
$ws.Range('C57').Value = 'RT @porrachatu: "você deveria ser mais legal" https://t.co/BFD7wsLkm0'
$ws.Range('C58').Value = 'RT @michaeldweiss: I can''t even get my child to put her raincoat on in exchange for Teletubby time. Russia''s youth is enterprising. https:/…'
$ws.Range('C59').Value = 'RT @ArmaTorlk: Pour gagner votre précommande RT + follow me et @armateam :) ! Tirage au sort dimanche 19 à 20h ! #Torlk2017 https://t.co/1X…'
$ws.Range('C60').Value = 'RT @ehdaora: essa eh pra se apaixona https://t.co/N5MxOFAb03'
$ws.Range('C61').Value = 'RT @Noelia_Mansilla: Ojala que nunca necesite uno! @Belu_Mansilla https://t.co/beD72KTmoR'
$ws.Range('C62').Value = 'RT @TeamBangtanCL: Info ??(!!!) 
¡Big Hit ha respondido ante las amenazas hacia Jimin! 
#ArmysWillProtectJimin https://t.co/Q1iIAnQ0FF'
$ws.Range('C63').Value = 'Horario de Noticias CNTP en cntpradio. https://t.co/1M2O36vncn https://t.co/DxtrQwETVI'
$ws.Range('C64').Value = 'RT @BTS_twt: #노츄
#커밍순 https://t.co/lRJy5FkevO'
$ws.Range('C65').Value = 'RT @paulinaromo69: #NewProfilePic https://t.co/nRT2H9l4Z1'
$ws.Range('C66').Value = 'RT @TheCCCompanies: Living Life Like A G Mixtape Dropping April 28th ??????@Art_Gretzky PROD. BY @DHoodNational @SwaggBBeatz ?????? #indiemusic #…'
$ws.Range('C67').Value = '@BlogdenWelttag sagt mal, kann ich das eigentlich auch noch bearbeiten? Also die Twitter- und Facebookzeile zum Beispiel? ??'
$ws.Range('C68').Value = 'RT @Uber_Pix: Welcome little turtle https://t.co/Q32PdGDUul'
$ws.Range('C69').Value = 'RT @abcfree56: #歳納京子生誕祭2017 
#歳納京子
#ゆるゆり
#RTした人全員フォローする 
#いいねした人全員フォロー 
京子おめでとう????
トメィトゥトメィトゥ?? https://t.co/lR1hdvJL3c'
$ws.Range('C70').Value = 'RT @AmanatUllah23: भारत है, तो हम हैं। भारत की प्रगति में ही हर भारतीय की उन्नति है https://t.co/2Te5LSPsCV'
$ws.Range('C71').Value = 'This is my favourite Mariah song idk why she never liked this song tbh https://t.co/g3HOSLQyqB'
$ws.Range('C72').Value = 'RT @skinhub: ⭐️ Butterfly Knife | Ultraviolet Giveaway
* RT &amp; Follow
* Reply w/ Skinhub User ID
Winner in 24 Hours! https://t.co/lYNnEezz…'
$ws.Range('C73').Value = 'RT @LicCarlosSosa: Con ayuda de Dip. Fermín Trujillo, se apoyo al niño Angel Miguel Buitimea Pérez, talento Sonorense.
#NuevaAlianzaimpulsa…'
$ws.Range('C74').Value = 'RT @itsyourgamerguy: Had some incredible support for this photo over on Instagram! Great addition to #ForzaHorizon3 @ForzaMotorsport @Xbox…'
$ws.Range('C75').Value = 'RT @hhuyrrtty04: #사설토토사이트추천
#사다리사이트추천
메-이-저-놀-이-터
♏안전최고♏
♐♐♐♐♐
❄https://t.co/jDR4J3anrs❄☜바로가기
Ⓜ매일e벤~Ⓜ
♌♌♌♌♌
☃편안히즐기세요☃
☑hello☑ https:…'
$ws.Range('C76').Value = 'RT @DeOlhoBBBrasil: Amoreees!
Vamos lutar pelo trio #MallyMar e contra:
Falso moralismo
Arrogância
Ganância
Inveja
Calunia
Difamação
#ForaD…'
$ws.Range('C77').Value = 'RT @FAVELADOANTARES: Favelado também é gente!
NÃO DEIXEM ESSE VIDEO MORRER! https://t.co/7QNiUWDS6n'
$ws.Range('C78').Value = 'Omg @JackMorlenMusic voice is just???? https://t.co/mNvlffw5Qu'
$ws.Range('C79').Value = 'RT @awecoupIes: this is the perfect representation of me in haunted houses https://t.co/0C1E7PFzAD'
$ws.Range('C80').Value = 'RT @DucaVisko: Dobro, ima li vođa bar dva validna potpisa od onih 650? Jedan demant za drugim... https://t.co/b3lJQbOWFc'
$ws.Range('C81').Value = 'Share your thoughts: https://t.co/pkub4F5nO3 https://t.co/0FNc6D1ios'
$ws.Range('C82').Value = 'Cuando Mami y Papi descubren un feature nuevo en el celular ?? https://t.co/3sXNqwOLig'
$ws.Range('C83').Value = 'RT @badman_sean: how''s it goinggggggg https://t.co/WseFXKlTil'
$ws.Range('C84').Value = 'Héctor Suárez va a parar al hospital tras sufrir una fuerte caída https://t.co/9kRoKRPQ2D'
$ws.Range('C85').Value = '@narendramodi https://t.co/569xmO4gcB'
$ws.Range('C86').Value = 'RT @ErkanPusmaz: https://t.co/zrnUjqnSLu'
$ws.Range('C87').Value = '#GirlPower #WomenLead https://t.co/nGD76ESDlp'
$ws.Range('C88').Value = 'RT @OulivierJirou: Lavezzi cet arnaqueur ?? https://t.co/MtmajyqkaE'
$ws.Range('C89').Value = 'RT @FirstTake: "LeBron has done so much for the game... He''s earned the opportunity to take a rest." - @kobebryant https://t.co/frQkcalV3N'
$ws.Range('C90').Value = 'RT @amam_990: @ajmi604 
•??
??
??
#الامير_العجمى_100k_محب
شاهدين لـكـ بـالخير يارب
•??
??
??
#حساب_أفتخر_فيه
#حساب_ملكي
????????????
@ajmi604 
•??
??
#تو…'
$ws.Range('C91').Value = 'RT @mjulio777pr: I know y''all remember this one. Wisin''s verse was ?? https://t.co/7H4UQE8eU6'
$ws.Range('C92').Value = 'RT @justjamiie: https://t.co/0shbb8bNJe'
$ws.Range('C93').Value = 'RT @hanxine: Slytherin don''t deserve this https://t.co/VO3MBW2v3v'
$ws.Range('C94').Value = 'RT @3zuwan: Lamborghini Huracan VS 14'' CBR 1000RR. Rilek je Lamborghini ni kena tapau. https://t.co/SstzsHdpzN'
$ws.Range('C95').Value = 'RT @waterparks: HI MOM. WE GOT OUR FIRST MAGAZINE COVER. THIS IS THE COOLEST THING IN THE WORLD. DAMN. https://t.co/tTBFAZuUUj https://t.co…'
$ws.Range('C96').Value = 'RT @JaMir_Russell: Pure evil lol https://t.co/Is3BIz7rMw'
$ws.Range('C97').Value = 'RT @caacosta1962: Uds. saben porqué @Lenin no quiso debatir?
No?
Yo si se
Porque no lo dejaron contar cachos de $10,000 https://t.co/rsx…'
$ws.Range('C98').Value = 'RT @BarstoolUA: Monday morning blues? This should help.
https://t.co/LvkcCJgBKF'
$ws.Range('C99').Value = 'RT @heyifeellike: when you''re dead inside but you still wanna have a good time. https://t.co/cetDuxQlSi'
$ws.Range('C100').Value = 'RT @MacCocktail: "When I sing, trouble can sit right on my shoulder and I don''t even notice." 
― Sarah Vaughan (born this day, March 27, 19…'
$ws.Range('C101').Value = 'RT @Stalinonyou: "Get in loser, we''re going shooting" https://t.co/DZVF5NGTVX'
